$wb = $excel.ActiveWorkbook

# --- 1. OpCodes sheet: fix CMP row's L8 cell. It previously held the stray
#        "r3" label; it should read "rDest" like the other ALU rows. Changing
#        it drops the now-unused "r3" shared string and the string table
#        renumbers itself automatically.
$opcodes = $wb.Worksheets.Item("OpCodes")
$opcodes.Range("L8").Value = "rDest"
$opcodes.Range("C6").Select()

# --- 2. Rename Sheet3 -> Resources and populate it with the new ALU
#        resource-usage table.
$res = $wb.Worksheets.Item("Sheet3")
$res.Name = "Resources"

$res.Range("A1").Value = "R32V2020 Resources"

$res.Range("A2").Value = "Module"
$res.Range("B2").Value = "L.E.s"
$res.Range("C2").Value = "Registers"
$res.Range("D2").Value = "Memory"
$res.Range("E2").Value = "Multipliers"
$res.Range("F2").Value = "PLLs"

$res.Range("A3").Value = "OpCodeDecoder"
$res.Range("B3").Value = 54
$res.Range("C3").Value = 0
$res.Range("D3").Value = 0
$res.Range("E3").Value = 0
$res.Range("F3").Value = "--"

$res.Range("A4").Value = "ALU"
$res.Range("B4").Value = 189
$res.Range("C4").Value = 0
$res.Range("D4").Value = 0
$res.Range("E4").Value = 0
$res.Range("F4").Value = 0

$res.Range("A5").Value = "Register File"
$res.Range("B5").Value = 594
$res.Range("C5").Value = 416
$res.Range("D5").Value = 0
$res.Range("E5").Value = 0
$res.Range("F5").Value = 0

$res.Range("A6").Value = "MUX_16x32"
$res.Range("B6").Value = 321
$res.Range("C6").Value = 416
$res.Range("D6").Value = 0
$res.Range("E6").Value = 0
$res.Range("F6").Value = 0

$res.Range("B2:F9").HorizontalAlignment = -4108

$res.Range("B11").Formula = "=SUM(B3:B10)"
$res.Range("C11").Formula = "=SUM(C3:C10)"
$res.Range("D11").Formula = "=SUM(D3:D10)"
$res.Range("E11").Formula = "=SUM(E3:E10)"
$res.Range("F11").Formula = "=SUM(F3:F10)"

$res.Columns.Item(1).ColumnWidth = 18.44
$res.Columns.Item(2).ColumnWidth = 5.18
$res.Columns.Item(3).ColumnWidth = 8.85
$res.Columns.Item(4).ColumnWidth = 8.6
$res.Columns.Item(6).ColumnWidth = 5.06

$res.Range("A11").Select()
